# Cambio menu color azul
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: move the "Ajustar php word..." task from column B (EN PROCESO)
# to column A (PENDIENTES) -- keep its existing cell style.
$val37 = $ws.Range("B37").Value2
$ws.Range("B37").Value = $null
$ws.Range("A37").Value = $val37

# Row 39: add new task "maquetear el proyecto parte visual" to column B (EN PROCESO)
$ws.Range("B39").Value = "maquetear el proyecto parte visual"

# Row 40: add new task "solucionar conflictos excel facturación" to column A (PENDIENTES)
$ws.Range("A40").Value = "solucionar conflictos excel facturación"

# Update the active selection to reflect the newly added rows further
# down the sheet (mirrors the saved cursor position in the sheet view).
$ws.Range("A41").Select()
